$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "10×73=",
    "16×92=",
    "48×13=",
    "73×48=",
    "22×59=",
    "61×86=",
    "90×51=",
    "85×65=",
    "33×54=",
    "45×44=",
    "29×84=",
    "76×56=",
    "100×49=",
    "79×96=",
    "25×19=",
    "96×91=",
    "71×37=",
    "48×53=",
    "25×54=",
    "26×84=",
    "34×79=",
    "24×44=",
    "95×51=",
    "56×84=",
    "24×51=",
    "93×36=",
    "85×72=",
    "24×75=",
    "83×23=",
    "77×91=",
    "39×52=",
    "43×20=",
    "38×60=",
    "37×83=",
    "44×15=",
    "86×13=",
    "54×30=",
    "95×14=",
    "26×53=",
    "80×31=",
    "94×88=",
    "93×76=",
    "72×36=",
    "100×100=",
    "81×29=",
    "78×83=",
    "29×87=",
    "63×83=",
    "48×88=",
    "64×31=",
    "55×67=",
    "88×19=",
    "84×20=",
    "24×62=",
    "53×44=",
    "29×31=",
    "23×48=",
    "34×76=",
    "65×78=",
    "69×82=",
    "22×85=",
    "31×48=",
    "59×34=",
    "46×87=",
    "86×20=",
    "60×85=",
    "88×50=",
    "10×18=",
    "73×70=",
    "57×37=",
    "79×78=",
    "77×85=",
    "59×48=",
    "55×38=",
    "75×90=",
    "38×43=",
    "71×55=",
    "85×57=",
    "82×63=",
    "30×38=",
    "31×23=",
    "24×39=",
    "53×75=",
    "81×18=",
    "32×46=",
    "36×20=",
    "96×78=",
    "28×32=",
    "27×86=",
    "91×41=",
    "14×41=",
    "50×29=",
    "12×59=",
    "95×25=",
    "97×24=",
    "80×11=",
    "81×89=",
    "56×11=",
    "22×57=",
    "57×52="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output ("Updated cells: " + $idx)
